$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - account holder first name
$ws.Range("C2").Value = "Hartmut"

# Row 3 - card number (keep as text) and account holder last name
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5 - opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 15.07.2025"

# Row 6 - first transaction
$ws.Range("B6").Value = "16.07."
$ws.Range("C6").Value = "17.07."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-57729601"
$ws.Range("E6").Value = "57,12-"

# Row 7 - second transaction
$ws.Range("B7").Value = "17.07."
$ws.Range("C7").Value = "18.07."
$ws.Range("D7").Value = "PAYPAL ROGFRD"
$ws.Range("E7").Value = "46,55-"

# Row 8 - third transaction
$ws.Range("B8").Value = "21.07."
$ws.Range("C8").Value = "22.07."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "24,63-"

# Row 12 - closing balance date and amount
$ws.Range("D12").Value = "KONTOSTAND AM 26.07.2025"
$ws.Range("E12").Value = "128,30-"

# Row 13 - next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 03.08.2025"
